$d = $word.ActiveDocument

# 1. Insert "import os" as a new paragraph right after "import socket"
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$d.Paragraphs.Item(2).Range.Text = "import os"

# 2. Rename the parameter in the function signature: save_path -> save_dir
$d.Content.Find.Execute(
    "def receive_file(save_path, port):", $true, $false, $false, $false,
    $false, $true, 1, $false, "def receive_file(save_dir, port):", 2
) | Out-Null

# 3. After "    sock.bind(('', port))" insert a block that receives the
#    filename and builds the full save_path, before the "with open(...)" line.
$bindIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "    sock.bind(('', port))`r") {
        $bindIdx = $i
        break
    }
}

$newLines = @(
    "",
    "    # Receive the filename first",
    "    filename, _ = sock.recvfrom(CHUNK_SIZE)",
    "    filename = filename.decode()",
    "",
    "    # Full path to save the file",
    "    save_path = os.path.join(save_dir, filename)",
    ""
)

$idx = $bindIdx
foreach ($line in $newLines) {
    $d.Paragraphs.Item($idx).Range.InsertParagraphAfter()
    $idx = $idx + 1
    if ($line -ne "") {
        $d.Paragraphs.Item($idx).Range.Text = $line
    }
}
